$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new value would otherwise be
# auto-converted to a number by Excel (values like "1.000", "330.25", etc.),
# matching the source data which stores these figures as plain text.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.134.18"
$ws.Range("E2").Value = "  +5.74%  "
$ws.Range("D3").Value = "1.920.99"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.64%  "
$ws.Range("D5").Value = "330.25"
$ws.Range("E5").Value = "  +4.72%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("D7").Value = "0.5220"
$ws.Range("E7").Value = "  +3.06%  "
$ws.Range("D8").Value = "0.4084"
$ws.Range("E8").Value = "  +4.43%  "
$ws.Range("D9").Value = "0.08513"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("E10").Value = "  +1.96%  "
$ws.Range("D11").Value = "42.80"
$ws.Range("E11").Value = "  +2.47%  "
$ws.Range("D12").Value = "22.39"
$ws.Range("E12").Value = "  +9.68%  "
$ws.Range("D13").Value = "6.429"
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("D14").Value = "1.916.65"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "7.410"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "95.74"
$ws.Range("E17").Value = "  +5.07%  "
$ws.Range("D18").Value = "0.00001114"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "0.06692"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").Value = "  +3.83%  "
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "6.016"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("D23").Value = "30.131.68"
$ws.Range("E23").Value = "  +5.66%  "
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("D25").Value = "2.209"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").Value = "2.146.05"
$ws.Range("E26").Value = "  +2.88%  "
$ws.Range("D27").Value = "21.13"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").Value = "159.86"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").Value = "2.449"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "129.46"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("D31").Value = "1.084"
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("D33").Value = "6.059"
$ws.Range("E33").Value = "  +5.71%  "
$ws.Range("D34").Value = "3.630"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("D36").Value = "0.06628"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").Value = "0.2216"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").Value = "1.239"
$ws.Range("E38").Value = "  +4.90%  "
$ws.Range("D39").Value = "5.196"
$ws.Range("E39").Value = "  +3.48%  "
$ws.Range("D40").Value = "8.907"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "0.6554"
$ws.Range("E41").Value = "  +2.73%  "
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").Value = "11.63"
$ws.Range("E43").Value = "  +4.86%  "
$ws.Range("D44").Value = "0.6166"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("D45").Value = "13.24"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").Value = "2.084"
$ws.Range("E47").Value = "  +3.91%  "
$ws.Range("D48").Value = "1.250"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("D49").Value = "124.66"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").Value = "1.175"
$ws.Range("E50").Value = "  +11.71%  "
$ws.Range("D51").Value = "79.80"
$ws.Range("E51").Value = "  +4.14%  "
